$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column headers: K = multiTexture, L = damage
$ws.Range("K1").Value = "multiTexture"
$ws.Range("L1").Value = "damage"

# multiTexture values for rows 2-12 (air, grass, dirt, stone, flower, tree x6)
$multiTexture = @(0,1,1,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $multiTexture.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $multiTexture[$i]
    $ws.Cells.Item($row, 12).Value = 0
}

# New row 13: spike block
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = "spike"
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 1

$ws.Range("L13").Select()
